$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107 (shifts existing rows 107-166 down to 108-167)
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(107, 1).Value = 4
$ws.Cells.Item(107, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(107, 3).Value = "Los Lagos"
$ws.Cells.Item(107, 4).Value = 45205
$ws.Cells.Item(107, 5).Value = 10
$ws.Cells.Item(107, 6).Value = 100112026
$ws.Cells.Item(107, 7).Value = "Haba"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 120
$ws.Cells.Item(107, 11).Value = 17000
$ws.Cells.Item(107, 12).Value = 17000
$ws.Cells.Item(107, 13).Value = 17000
$ws.Cells.Item(107, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(107, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(107, 16).Value = 680
$ws.Cells.Item(107, 17).Value = 25
$ws.Cells.Item(107, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the other date cells in column D
$ws.Cells.Item(107, 4).NumberFormat = $ws.Cells.Item(106, 4).NumberFormat
